$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2483443708609271
$ws.Range("C2").Value = 0.4503311258278146
$ws.Range("J2").Value = 0.02980132450331126
$ws.Range("P2").Value = 0.152317880794702
$ws.Range("S2").Value = 0.119205298013245
$ws.Range("B3").Value = 0.01379310344827586
$ws.Range("C3").Value = 0.04137931034482759
$ws.Range("J3").Value = 0.03448275862068965
$ws.Range("P3").Value = 0.6758620689655173
$ws.Range("S3").Value = 0.2344827586206897
$ws.Range("P4").Value = 0.5757575757575758
$ws.Range("S4").Value = 0.4242424242424243
$ws.Range("B6").Value = 0.05829596412556054
$ws.Range("D6").Value = 0.008968609865470852
$ws.Range("F6").Value = 0.1031390134529148
$ws.Range("J6").Value = 0.2331838565022422
$ws.Range("O6").Value = 0.02690582959641256
$ws.Range("Q6").Value = 0.1300448430493273
$ws.Range("R6").Value = 0.05381165919282511
$ws.Range("S6").Value = 0.3856502242152466
$ws.Range("B7").Value = 0.07602339181286549
$ws.Range("D7").Value = 0.01754385964912281
$ws.Range("F7").Value = 0.04678362573099415
$ws.Range("J7").Value = 0.1578947368421053
$ws.Range("O7").Value = 0.03508771929824561
$ws.Range("Q7").Value = 0.1754385964912281
$ws.Range("R7").Value = 0.1052631578947368
$ws.Range("S7").Value = 0.3859649122807017
$ws.Range("B8").Value = 0.09959349593495935
$ws.Range("D8").Value = 0.008130081300813009
$ws.Range("F8").Value = 0.06097560975609756
$ws.Range("J8").Value = 0.1199186991869919
$ws.Range("O8").Value = 0.01829268292682927
$ws.Range("Q8").Value = 0.1707317073170732
$ws.Range("R8").Value = 0.07926829268292683
$ws.Range("S8").Value = 0.443089430894309
$ws.Range("B9").Value = 0.08433734939759036
$ws.Range("D9").Value = 0.006024096385542169
$ws.Range("F9").Value = 0.06024096385542169
$ws.Range("J9").Value = 0.1204819277108434
$ws.Range("O9").Value = 0.02409638554216868
$ws.Range("Q9").Value = 0.2048192771084337
$ws.Range("R9").Value = 0.1144578313253012
$ws.Range("S9").Value = 0.3855421686746988
$ws.Range("B10").Value = 0.1030684500393391
$ws.Range("D10").Value = 0.01888276947285602
$ws.Range("F10").Value = 0.06530291109362707
$ws.Range("J10").Value = 0.1313926042486231
$ws.Range("O10").Value = 0.01730920535011802
$ws.Range("Q10").Value = 0.2265932336742722
$ws.Range("R10").Value = 0.08733280881195908
$ws.Range("S10").Value = 0.3501180173092053
$ws.Range("G11").Value = 0.155893536121673
$ws.Range("J11").Value = 0.1026615969581749
$ws.Range("K11").Value = 0.1977186311787072
$ws.Range("L11").Value = 0.5285171102661597
$ws.Range("S11").Value = 0.01520912547528517
$ws.Range("G12").Value = 0.7132867132867133
$ws.Range("J12").Value = 0.2307692307692308
$ws.Range("K12").Value = 0.006993006993006993
$ws.Range("L12").Value = 0.02797202797202797
$ws.Range("S12").Value = 0.02097902097902098
$ws.Range("F15").Value = 0.02678571428571428
$ws.Range("H15").Value = 0.1919642857142857
$ws.Range("I15").Value = 0.05803571428571429
$ws.Range("J15").Value = 0.3125
$ws.Range("K15").Value = 0.04910714285714286
$ws.Range("M15").Value = 0.008928571428571428
$ws.Range("O15").Value = 0.05803571428571429
$ws.Range("S15").Value = 0.2946428571428572
$ws.Range("F17").Value = 0.01298701298701299
$ws.Range("H17").Value = 0.2207792207792208
$ws.Range("I17").Value = 0.06060606060606061
$ws.Range("J17").Value = 0.461038961038961
$ws.Range("K17").Value = 0.08008658008658008
$ws.Range("M17").Value = 0.01948051948051948
$ws.Range("O17").Value = 0.05844155844155844
$ws.Range("S17").Value = 0.08658008658008658
$ws.Range("F18").Value = 0.02030456852791878
$ws.Range("H18").Value = 0.1878172588832487
$ws.Range("I18").Value = 0.06091370558375635
$ws.Range("J18").Value = 0.4568527918781726
$ws.Range("K18").Value = 0.1116751269035533
$ws.Range("M18").Value = 0.01015228426395939
$ws.Range("N18").Value = 0.005076142131979695
$ws.Range("O18").Value = 0.06091370558375635
$ws.Range("S18").Value = 0.08629441624365482
$ws.Range("F19").Value = 0.01515151515151515
$ws.Range("H19").Value = 0.2352472089314195
$ws.Range("I19").Value = 0.07974481658692185
$ws.Range("J19").Value = 0.3580542264752791
$ws.Range("K19").Value = 0.09409888357256778
$ws.Range("M19").Value = 0.02392344497607655
$ws.Range("O19").Value = 0.07336523125996811
$ws.Range("S19").Value = 0.120414673046252
